$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-shuffled duty roster: update column B (name_duty) values row by row.
$ws.Range("B2").Value  = "Cox Matthew Jonah"
$ws.Range("B3").Value  = "Hansen Jakob U"
$ws.Range("B4").Value  = "石井海成"
$ws.Range("B5").Value  = "Nicholas Tristan Aryasatyo"
$ws.Range("B6").Value  = "小溝賢"
$ws.Range("B7").Value  = "小野文哉"
$ws.Range("B8").Value  = "渡部魁"
$ws.Range("B9").Value  = "崎谷航平"
$ws.Range("B10").Value = "三神佳誠"
$ws.Range("B11").Value = "氏家琉貴"
$ws.Range("B12").Value = "羽賀尚生"
$ws.Range("B13").Value = ""
$ws.Range("B14").Value = "島田実"
$ws.Range("B15").Value = "足立耕平"
$ws.Range("B16").Value = "遠藤隼人"
$ws.Range("B17").Value = "富澤天音"
$ws.Range("B18").Value = "神山修造"
$ws.Range("B19").Value = "川田涼介"
$ws.Range("B20").Value = "豊島亮"
$ws.Range("B21").Value = "兒島大志郎"
$ws.Range("B22").Value = "山口玲"
$ws.Range("B23").Value = "日高泰聖"
$ws.Range("B24").Value = "志塚惇希"
$ws.Range("B25").Value = "白岩詩佑介"
$ws.Range("B26").Value = "Cox Matthew Jonah"
$ws.Range("B27").Value = "Hansen Jakob U"
$ws.Range("B28").Value = "石井海成"
$ws.Range("B29").Value = "Nicholas Tristan Aryasatyo"
$ws.Range("B30").Value = "小溝賢"
$ws.Range("B31").Value = "小野文哉"
$ws.Range("B32").Value = "渡部魁"

# B29 previously carried a one-off explicit-black Arial 10pt font; normalize it
# back to the same font the rest of the column uses so the stray font/style
# definition drops out.
$ws.Range("B29").ClearFormats()
$ws.Range("B29").Font.Name = "Arial"
$ws.Range("B29").Font.Size = 10

# Restore the B2:B32 selection (with B2 as the active cell) that was in place
# when the sheet was saved.
$ws.Range("B2:B32").Select() | Out-Null
